$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting rows 30..107 down to 31..108.
$ws.Rows.Item(30).Insert()

# Fill in the new record (weekly Choclo / Hortaliza price entry) that was
# inserted at row 30.
$ws.Range("A30").Value = 7
$ws.Range("B30").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C30").Value = "Ñuble"
$ws.Range("D30").Value = 44533
$ws.Range("E30").Value = 16
$ws.Range("F30").Value = 100112024
$ws.Range("G30").Value = "Choclo"
$ws.Range("H30").Value = "Dulce o Americano"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 60
$ws.Range("K30").Value = 16000
$ws.Range("L30").Value = 17000
$ws.Range("M30").Value = 16500
$ws.Range("N30").Value = "$/malla 60 unidades"
$ws.Range("O30").Value = "Provincia de Limarí"
$ws.Range("P30").Value = 275
$ws.Range("Q30").Value = 60
$ws.Range("R30").Value = "Hortaliza"
